$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3 header info
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay TEXT (not be coerced to a
# number) - format the cell as Text first, assign the value, then
# paste the original cell's formatting back on top so the style index
# (s="8") is preserved (only the value/type changes, matching the diff).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 03.02.2025"

# Row 6
$ws.Range("B6").Value = "04.02."
$ws.Range("C6").Value = "05.02."
$ws.Range("D6").Value = "KARTENZ./04.02 LIDL RO"
$ws.Range("E6").Value = "128,60-"

# Row 7
$ws.Range("B7").Value = "05.02."
$ws.Range("C7").Value = "06.02."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "25,31-"

# Row 8
$ws.Range("B8").Value = "06.02."
$ws.Range("C8").Value = "07.02."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-30733188"
$ws.Range("E8").Value = "56,45-"

# Row 9 - previously a blank placeholder row; now filled with a new
# transaction. Set the values first, then paste row 8's formatting over
# it (B9:D9 -> s="8" like B8:D8, E9 -> s="17" like E8) so the style
# indexes match what a real new transaction row looks like.
$ws.Range("B9").Value = "09.02."
$ws.Range("C9").Value = "10.02."
$ws.Range("D9").Value = "BURGER KING Geithain"
$ws.Range("E9").Value = "30,21-"
$ws.Range("B8:E8").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)

# Row 10 - same treatment, another previously-blank placeholder row.
$ws.Range("B10").Value = "10.02."
$ws.Range("C10").Value = "11.02."
$ws.Range("D10").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E10").Value = "46,79-"
$ws.Range("B9:E9").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 13.02.2025"
$ws.Range("E12").Value = "287,36-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 19.02.2025"
